$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1102
$ws1.Range("F8").Value = 392
$ws1.Range("F14").Value = 12630
$ws1.Range("F15").Value = 5206
$ws1.Range("F16").Value = 5518

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 34

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1102
$ws4.Range("F9").Value = 392
$ws4.Range("F15").Value = 12630
$ws4.Range("F16").Value = 34
$ws4.Range("F18").Value = 5206
$ws4.Range("F19").Value = 5518
